$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 910.2727
$ws.Range("I96").Value = 575.1429000000001
$ws.Range("J96").Value = 1496.75
$ws.Range("K96").Value = 1725.4287
$ws.Range("L96").Value = 4490.25
$ws.Range("M96").Value = -352.4287000000002
$ws.Range("N96").Value = -7236.25

$ws.Range("H135").Value = 867.7442
$ws.Range("I135").Value = 514.8333
$ws.Range("J135").Value = 2682.7144
$ws.Range("K135").Value = 4633.4997
$ws.Range("L135").Value = 24144.4296
$ws.Range("M135").Value = -2098.4997
$ws.Range("N135").Value = -29214.4296

$ws.Range("H137").Value = 7343.1226
$ws.Range("I137").Value = 5801.9
$ws.Range("J137").Value = 8406.034
$ws.Range("K137").Value = 17405.7
$ws.Range("L137").Value = 25218.102
$ws.Range("M137").Value = -14855.7
$ws.Range("N137").Value = -30318.102

$ws.Range("H138").Value = 2685.2354
$ws.Range("I138").Value = 1372.8572
$ws.Range("J138").Value = 4805.231
$ws.Range("K138").Value = 4118.571599999999
$ws.Range("L138").Value = 14415.693
$ws.Range("M138").Value = 1021.428400000001
$ws.Range("N138").Value = -24695.693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5409.357
$ws.Range("I74").Value = 8433.846
$ws.Range("J74").Value = 4053.5518
$ws.Range("K74").Value = 8433.846
$ws.Range("L74").Value = 4053.5518
$ws.Range("M74").Value = -7559.846
$ws.Range("N74").Value = -5801.5518

$ws.Range("H77").Value = 5409.357
$ws.Range("I77").Value = 8433.846
$ws.Range("J77").Value = 4053.5518
$ws.Range("K77").Value = 42169.23
$ws.Range("L77").Value = 20267.759
$ws.Range("M77").Value = -37801.23
$ws.Range("N77").Value = -29003.759

$ws.Range("H110").Value = 1575.6875
$ws.Range("I110").Value = 1413.2
$ws.Range("J110").Value = 4013
$ws.Range("K110").Value = 1413.2
$ws.Range("L110").Value = 4013
$ws.Range("M110").Value = 631.8
$ws.Range("N110").Value = -8103

$ws.Range("H134").Value = 31960
$ws.Range("J134").Value = 31960
$ws.Range("L134").Value = 31960
$ws.Range("N134").Value = -42100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1064.0416
$ws.Range("I99").Value = 766.2
$ws.Range("J99").Value = 1560.4445
$ws.Range("K99").Value = 766.2
$ws.Range("L99").Value = 1560.4445
$ws.Range("M99").Value = 731.8
$ws.Range("N99").Value = -4556.4445

$ws.Range("H107").Value = 2258.2856
$ws.Range("I107").Value = 2048.75
$ws.Range("K107").Value = 2048.75
$ws.Range("M107").Value = -128.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 109.53333
$ws.Range("I7").Value = 50.333332
$ws.Range("J7").Value = 124.333336
$ws.Range("K7").Value = 50.333332
$ws.Range("L7").Value = 124.333336
$ws.Range("M7").Value = 62.666668
$ws.Range("N7").Value = -350.333336

$ws.Range("H58").Value = 871.2353000000001
$ws.Range("I58").Value = 1154.5416
$ws.Range("J58").Value = 619.4074000000001
$ws.Range("K58").Value = 1154.5416
$ws.Range("L58").Value = 619.4074000000001
$ws.Range("M58").Value = -951.5416
$ws.Range("N58").Value = -1025.4074

$ws.Range("H132").Value = 34488660
$ws.Range("I132").Value = 58832240
$ws.Range("J132").Value = 1918.6666
$ws.Range("K132").Value = 176496720
$ws.Range("L132").Value = 5755.9998
$ws.Range("M132").Value = -176494190
$ws.Range("N132").Value = -10815.9998

$ws.Range("H136").Value = 871.2353000000001
$ws.Range("I136").Value = 1154.5416
$ws.Range("J136").Value = 619.4074000000001
$ws.Range("K136").Value = 3463.6248
$ws.Range("L136").Value = 1858.2222
$ws.Range("M136").Value = -913.6248000000001
$ws.Range("N136").Value = -6958.2222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 294.93332
$ws.Range("I5").Value = 306.21054
$ws.Range("J5").Value = 275.45456
$ws.Range("K5").Value = 918.6316199999999
$ws.Range("L5").Value = 826.36368
$ws.Range("M5").Value = -806.6316199999999
$ws.Range("N5").Value = -1050.36368

$ws.Range("H68").Value = 988.875
$ws.Range("I68").Value = 601.4103
$ws.Range("K68").Value = 1804.2309
$ws.Range("M68").Value = -993.2309

$ws.Range("H71").Value = 988.875
$ws.Range("I71").Value = 601.4103
$ws.Range("K71").Value = 5412.6927
$ws.Range("M71").Value = -1356.6927

$ws.Range("H107").Value = 450.34042
$ws.Range("I107").Value = 196.45
$ws.Range("J107").Value = 1901.1428
$ws.Range("K107").Value = 589.3499999999999
$ws.Range("L107").Value = 5703.428400000001
$ws.Range("M107").Value = 1330.65
$ws.Range("N107").Value = -9543.428400000001

$ws.Range("H131").Value = 634.7041
$ws.Range("I131").Value = 210.19513
$ws.Range("J131").Value = 940.0526
$ws.Range("K131").Value = 630.58539
$ws.Range("L131").Value = 2820.1578
$ws.Range("M131").Value = 4409.41461
$ws.Range("N131").Value = -12900.1578

$ws.Range("H135").Value = 294.93332
$ws.Range("I135").Value = 306.21054
$ws.Range("J135").Value = 275.45456
$ws.Range("K135").Value = 2755.89486
$ws.Range("L135").Value = 2479.09104
$ws.Range("M135").Value = -220.8948599999999
$ws.Range("N135").Value = -7549.09104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 30033.875
$ws.Range("J64").Value = 30033.875
$ws.Range("L64").Value = 30033.875
$ws.Range("N64").Value = -30529.875

$ws.Range("H67").Value = 30033.875
$ws.Range("J67").Value = 30033.875
$ws.Range("L67").Value = 30033.875
$ws.Range("N67").Value = -31749.875

$ws.Range("H126").Value = 1792.2142
$ws.Range("I126").Value = 1673.875
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 5021.625
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -2551.625
$ws.Range("N126").Value = -10790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3349.5833
$ws.Range("I136").Value = 1270.0344
$ws.Range("J136").Value = 6523.6313
$ws.Range("K136").Value = 3810.1032
$ws.Range("L136").Value = 19570.8939
$ws.Range("M136").Value = -1260.1032
$ws.Range("N136").Value = -24670.8939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5744069.5
$ws.Range("I122").Value = 7179586.5
$ws.Range("J122").Value = 2000.7142
$ws.Range("K122").Value = 21538759.5
$ws.Range("L122").Value = 6002.142599999999
$ws.Range("M122").Value = -21536309.5
$ws.Range("N122").Value = -10902.1426

$ws.Range("H136").Value = 2584.4924
$ws.Range("I136").Value = 3780.4856
$ws.Range("J136").Value = 1276.375
$ws.Range("K136").Value = 11341.4568
$ws.Range("L136").Value = 3829.125
$ws.Range("M136").Value = -8791.4568
$ws.Range("N136").Value = -8929.125
